# Actualizacion 9 marzo 2025
#
# Fix capitalization of the "Maestría" education entry title so the
# shared-string text matches proper title case, and update the
# active-cell selection left in the sheet view (as last saved by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the capitalization of the master's degree title (row 3, col A):
#   "Maestría en Igualdad de género en ámbito público y privado "
#   -> "Maestría en Igualdad de Género en Ámbito Público y Privado "
$ws.Cells.Item(3, 1).Value = "Maestría en Igualdad de Género en Ámbito Público y Privado "

# Update the sheet's remembered selection (cursor left on C10 on last save).
[void]$ws.Range("C10").Select()
